{"js": "const body = context.document.body;\n\n// Update the title / date line.\nconst titleHits = body.search(\"2024-06-06 Thursday\", { matchCase: true });\ntitleHits.load(\"items\");\nawait context.sync();\nif (titleHits.items.length > 0) {\n  titleHits.items[0].insertText(\"2024-06-07 Friday\", \"Replace\");\n}\n\n// Update the division problems laid out in the table. The table also\n// contains blank answer rows interleaved with the problem rows, so address\n// cells by row/column position (two of the original values repeat, so a\n// plain text search-and-replace would be ambiguous).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst rowsOld = [\n  [\"35\u00f73=\", \"32\u00f78=\", \"45\u00f74=\", \"55\u00f76=\", \"59\u00f78=\"],\n  [\"92\u00f73=\", \"78\u00f73=\", \"97\u00f72=\", \"19\u00f77=\", \"46\u00f72=\"],\n  [\"70\u00f72=\", \"32\u00f73=\", \"77\u00f74=\", \"53\u00f78=\", \"96\u00f77=\"],\n  [\"35\u00f72=\", \"41\u00f79=\", \"88\u00f73=\", \"74\u00f77=\", \"91\u00f72=\"],\n  [\"20\u00f78=\", \"17\u00f75=\", \"33\u00f79=\", \"32\u00f78=\", \"45\u00f74=\"],\n];\nconst rowsNew = [\n  [\"83\u00f73=\", \"92\u00f75=\", \"30\u00f74=\", \"10\u00f77=\", \"42\u00f76=\"],\n  [\"66\u00f73=\", \"21\u00f74=\", \"73\u00f77=\", \"93\u00f79=\", \"97\u00f75=\"],\n  [\"59\u00f73=\", \"11\u00f72=\", \"87\u00f72=\", \"95\u00f77=\", \"98\u00f79=\"],\n  [\"24\u00f75=\", \"91\u00f78=\", \"17\u00f78=\", \"55\u00f75=\", \"41\u00f72=\"],\n  [\"11\u00f78=\", \"77\u00f72=\", \"63\u00f74=\", \"12\u00f77=\", \"44\u00f74=\"],\n];\n\n// 0-based table row indexes that actually hold the division problems (the\n// rows in between are blank answer rows).\nconst problemRowIndexes = [0, 4, 8, 12, 16];\n\nconst cells = [];\nfor (let r = 0; r < problemRowIndexes.length; r++) {\n  for (let c = 0; c < 5; c++) {\n    const cell = table.getCell(problemRowIndexes[r], c);\n    cell.load(\"value\");\n    cells.push({ cell, expected: rowsOld[r][c], next: rowsNew[r][c] });\n  }\n}\nawait context.sync();\n\nfor (const { cell, expected, next } of cells) {\n  if (cell.value === expected) {\n    cell.value = next;\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the title / date line.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"2024-06-06 Thursday\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"2024-06-07 Friday\"\n$find.Execute([ref]$find.Text, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, 0, [ref]$false, [ref]$find.Replacement.Text, 2)\n\n# Update the division problems laid out in the table, row by row (the\n# table also contains blank answer rows interleaved with the problem\n# rows, so walk only the rows/cells that actually hold text).\n$table = $d.Tables.Item(1)\n\n$rowsOld = @(\n    @(\"35\u00f73=\", \"32\u00f78=\", \"45\u00f74=\", \"55\u00f76=\", \"59\u00f78=\"),\n    @(\"92\u00f73=\", \"78\u00f73=\", \"97\u00f72=\", \"19\u00f77=\", \"46\u00f72=\"),\n    @(\"70\u00f72=\", \"32\u00f73=\", \"77\u00f74=\", \"53\u00f78=\", \"96\u00f77=\"),\n    @(\"35\u00f72=\", \"41\u00f79=\", \"88\u00f73=\", \"74\u00f77=\", \"91\u00f72=\"),\n    @(\"20\u00f78=\", \"17\u00f75=\", \"33\u00f79=\", \"32\u00f78=\", \"45\u00f74=\")\n)\n$rowsNew = @(\n    @(\"83\u00f73=\", \"92\u00f75=\", \"30\u00f74=\", \"10\u00f77=\", \"42\u00f76=\"),\n    @(\"66\u00f73=\", \"21\u00f74=\", \"73\u00f77=\", \"93\u00f79=\", \"97\u00f75=\"),\n    @(\"59\u00f73=\", \"11\u00f72=\", \"87\u00f72=\", \"95\u00f77=\", \"98\u00f79=\"),\n    @(\"24\u00f75=\", \"91\u00f78=\", \"17\u00f78=\", \"55\u00f75=\", \"41\u00f72=\"),\n    @(\"11\u00f78=\", \"77\u00f72=\", \"63\u00f74=\", \"12\u00f77=\", \"44\u00f74=\")\n)\n\n$tableRowIndexes = @(1, 5, 9, 13, 17)\n\nfor ($r = 0; $r -lt $tableRowIndexes.Count; $r++) {\n    $tr = $table.Rows.Item($tableRowIndexes[$r])\n    for ($c = 1; $c -le 5; $c++) {\n        $cell = $tr.Cells.Item($c)\n        $expected = $rowsOld[$r][$c - 1]\n        $newVal = $rowsNew[$r][$c - 1]\n\n        $cellRange = $cell.Range\n        $cellRange.MoveEnd(1, -1) | Out-Null\n        if ($cellRange.Text -eq $expected) {\n            $cellRange.Text = $newVal\n        }\n    }\n}\n"}
